$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ParticipantsTab) query text (column B) was rewritten to a new Cypher
# query (adds diagnosis/genomic_info optional matches, re-derives samples via
# a second pass, sorts sample ids, and renames the ORDER BY keyword casing).
$newParticipantQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['JSON']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$ws.Range("B2").Value = $newParticipantQuery

# Update the view selection to match the saved workbook state (active cell
# moved from C5 to B4).
$ws.Range("B4").Select()
